$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 2 contents (the first exposure site - text updated)
$ws.Range("A2").Value = "Point Cook"
$ws.Range("B2").Value = "The Coffeeologist Cafe, 70/300 Point Cook Rd , Point Cook VIC 3030"
$ws.Range("C2").Value = "11:00am - 11:40am 8/2/2021"
$ws.Range("D2").Value = "Case attended venue"
$ws.Range("E2").Value = "new"

# Replace row 3 contents (second exposure site) - overwrite old data
$ws.Range("A3").Value = "Point Cook"
$ws.Range("B3").Value = "The Coffeeologist Cafe, 70/300 Point Cook Rd, Point Cook VIC 3030"
$ws.Range("C3").Value = "11:30am - 12:10pm 10/2/2021"
$ws.Range("D3").Value = "Case attended venue"
$ws.Range("E3").Value = "new"

# Remove old rows 4-6 (data no longer present)
$ws.Range("A4:E6").EntireRow.Delete()

# Auto-fit column widths to the new (shorter) content
$ws.Columns.Item(1).ColumnWidth = 9.19921875
$ws.Columns.Item(2).ColumnWidth = 54.53125
$ws.Columns.Item(3).ColumnWidth = 25.53125
$ws.Columns.Item(4).ColumnWidth = 17.265625

# Update selection to reflect active cell
$ws.Range("C2").Select()
